# Updates cryptos list values per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.753.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0637"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.53"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0782"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.640.62"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.25"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.856.89"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.771.07"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.42"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.57"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.94"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.23"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.120"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.82"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.49"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0486"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.39"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.546"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.106.29"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0155"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.59"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.77"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.800"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0110"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.02"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.68"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.33"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.39%  "
